# Apply scheduled-runner value updates to the Gilgamesh_Profits sheets.
# Workbook sheet tabs (ALC/ARM/.../WVR) correspond 1:1 to the physical
# sheet parts touched by the upstream diff.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2373
$ws.Range("J17").Value = 2373
$ws.Range("L17").Value = 7119
$ws.Range("N17").Value = -7455
$ws.Range("H68").Value = 59998
$ws.Range("J68").Value = 59998
$ws.Range("L68").Value = 59998
$ws.Range("N68").Value = -61496
$ws.Range("H71").Value = 59998
$ws.Range("J71").Value = 59998
$ws.Range("L71").Value = 179994
$ws.Range("N71").Value = -187482
$ws.Range("H132").Value = 4997.2
$ws.Range("I132").Value = 4997.2
$ws.Range("K132").Value = 14991.6
$ws.Range("M132").Value = -12461.6
$ws.Range("H138").Value = 275386.97
$ws.Range("I138").Value = 3446.4243
$ws.Range("J138").Value = 427489.28
$ws.Range("K138").Value = 10339.2729
$ws.Range("L138").Value = 1282467.84
$ws.Range("M138").Value = -5199.2729
$ws.Range("N138").Value = -1292747.84
$ws.Range("H141").Value = 3094.7896
$ws.Range("I141").Value = 2076.2307
$ws.Range("J141").Value = 5301.6665
$ws.Range("K141").Value = 6228.6921
$ws.Range("L141").Value = 15904.9995
$ws.Range("M141").Value = -1048.6921
$ws.Range("N141").Value = -26264.9995

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N35").ClearContents()
$ws.Range("H2").Value = 696.4583
$ws.Range("I2").Value = 663.1905
$ws.Range("J2").Value = 929.3333
$ws.Range("K2").Value = 663.1905
$ws.Range("L2").Value = 929.3333
$ws.Range("M2").Value = -550.1905
$ws.Range("N2").Value = -1155.3333
$ws.Range("H5").Value = 107.25
$ws.Range("I5").Value = 107.25
$ws.Range("K5").Value = 107.25
$ws.Range("M5").Value = 4.75
$ws.Range("H35").Value = 4742.2
$ws.Range("I35").Value = 4742.2
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 4742.2
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -4336.2
$ws.Range("H116").Value = 696.4583
$ws.Range("I116").Value = 663.1905
$ws.Range("J116").Value = 929.3333
$ws.Range("K116").Value = 663.1905
$ws.Range("L116").Value = 929.3333
$ws.Range("M116").Value = 1630.8095
$ws.Range("N116").Value = -5517.3333
$ws.Range("H139").Value = 76725.45
$ws.Range("J139").Value = 76725.45
$ws.Range("L139").Value = 76725.45
$ws.Range("N139").Value = -87005.45

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 696.4583
$ws.Range("I3").Value = 663.1905
$ws.Range("J3").Value = 929.3333
$ws.Range("K3").Value = 663.1905
$ws.Range("L3").Value = 929.3333
$ws.Range("M3").Value = -549.1905
$ws.Range("N3").Value = -1157.3333
$ws.Range("H4").Value = 107.25
$ws.Range("I4").Value = 107.25
$ws.Range("K4").Value = 107.25
$ws.Range("M4").Value = 7.75
$ws.Range("H20").Value = 26886526
$ws.Range("I20").Value = 34727570
$ws.Range("J20").Value = 2956
$ws.Range("K20").Value = 34727570
$ws.Range("L20").Value = 2956
$ws.Range("M20").Value = -34727323
$ws.Range("N20").Value = -3450
$ws.Range("H22").Value = 879.9231
$ws.Range("I22").Value = 886.4167
$ws.Range("K22").Value = 886.4167
$ws.Range("M22").Value = -713.4167
$ws.Range("H99").Value = 102047.6
$ws.Range("I99").Value = 125684.5
$ws.Range("K99").Value = 125684.5
$ws.Range("M99").Value = -124186.5
$ws.Range("H107").Value = 2263604.2
$ws.Range("I107").Value = 2748306
$ws.Range("K107").Value = 2748306
$ws.Range("M107").Value = -2746386
$ws.Range("H134").Value = 3531.75
$ws.Range("I134").Value = 2779.2354
$ws.Range("K134").Value = 8337.706200000001
$ws.Range("M134").Value = -5802.706200000001
$ws.Range("H141").Value = 74778.836
$ws.Range("J141").Value = 79734.60000000001
$ws.Range("L141").Value = 79734.60000000001
$ws.Range("N141").Value = -90094.60000000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2317.5715
$ws.Range("I94").Value = 1681.75
$ws.Range("J94").Value = 2571.9
$ws.Range("K94").Value = 1681.75
$ws.Range("L94").Value = 2571.9
$ws.Range("M94").Value = -1230.75
$ws.Range("N94").Value = -3473.9
$ws.Range("H132").Value = 11908116
$ws.Range("I132").Value = 14708909
$ws.Range("K132").Value = 44126727
$ws.Range("M132").Value = -44124197
$ws.Range("H134").Value = 2478.2415
$ws.Range("I134").Value = 2016.5652
$ws.Range("K134").Value = 6049.6956
$ws.Range("M134").Value = -3514.6956

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M17").ClearContents()
$ws.Range("H17").Value = 2500
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2500
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7500
$ws.Range("N17").Value = -7838
$ws.Range("H87").Value = 24323
$ws.Range("I87").Value = 24323
$ws.Range("K87").Value = 72969
$ws.Range("M87").Value = -71721
$ws.Range("H90").Value = 24323
$ws.Range("I90").Value = 24323
$ws.Range("K90").Value = 218907
$ws.Range("M90").Value = -212667
$ws.Range("H138").Value = 3974.4546
$ws.Range("I138").Value = 2302.5
$ws.Range("K138").Value = 6907.5
$ws.Range("M138").Value = -1767.5
$ws.Range("H139").Value = 7390.7617
$ws.Range("I139").Value = 9569.846
$ws.Range("J139").Value = 3849.75
$ws.Range("K139").Value = 28709.538
$ws.Range("L139").Value = 11549.25
$ws.Range("M139").Value = -23569.538
$ws.Range("N139").Value = -21829.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M35").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H35").Value = 15000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15596
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H70").Value = 16198383
$ws.Range("I70").Value = 23909184
$ws.Range("J70").Value = 5700.1
$ws.Range("K70").Value = 23909184
$ws.Range("L70").Value = 5700.1
$ws.Range("M70").Value = -23908914
$ws.Range("N70").Value = -6240.1
$ws.Range("H73").Value = 16198383
$ws.Range("I73").Value = 23909184
$ws.Range("J73").Value = 5700.1
$ws.Range("K73").Value = 23909184
$ws.Range("L73").Value = 5700.1
$ws.Range("M73").Value = -23908248
$ws.Range("N73").Value = -7572.1
$ws.Range("H97").Value = 3414.95
$ws.Range("I97").Value = 1190.2
$ws.Range("J97").Value = 5639.7
$ws.Range("K97").Value = 1190.2
$ws.Range("L97").Value = 5639.7
$ws.Range("M97").Value = -694.2
$ws.Range("N97").Value = -6631.7
$ws.Range("H126").Value = 8954.666999999999
$ws.Range("I126").Value = 2604.6667
$ws.Range("K126").Value = 7814.000100000001
$ws.Range("M126").Value = -5344.000100000001
$ws.Range("H132").Value = 2798.8333
$ws.Range("I132").Value = 2776.5186
$ws.Range("J132").Value = 2999.6667
$ws.Range("K132").Value = 8329.5558
$ws.Range("L132").Value = 8999.000100000001
$ws.Range("M132").Value = -5799.5558
$ws.Range("N132").Value = -14059.0001
$ws.Range("H133").Value = 113999.4
$ws.Range("J133").Value = 113999.4
$ws.Range("L133").Value = 113999.4
$ws.Range("N133").Value = -124119.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1492.6562
$ws.Range("I61").Value = 1412.9584
$ws.Range("J61").Value = 1731.75
$ws.Range("K61").Value = 1412.9584
$ws.Range("L61").Value = 1731.75
$ws.Range("M61").Value = -1210.9584
$ws.Range("N61").Value = -2135.75
$ws.Range("H100").Value = 1042849.3
$ws.Range("I100").Value = 1231717.4
$ws.Range("K100").Value = 1231717.4
$ws.Range("M100").Value = -1231176.4
$ws.Range("H113").Value = 1492.6562
$ws.Range("I113").Value = 1412.9584
$ws.Range("J113").Value = 1731.75
$ws.Range("K113").Value = 1412.9584
$ws.Range("L113").Value = 1731.75
$ws.Range("M113").Value = 757.0416
$ws.Range("N113").Value = -6071.75
$ws.Range("H122").Value = 3770.6428
$ws.Range("I122").Value = 4210.5557
$ws.Range("K122").Value = 12631.6671
$ws.Range("M122").Value = -10181.6671
$ws.Range("H138").Value = 82768.57000000001
$ws.Range("J138").Value = 82768.57000000001
$ws.Range("L138").Value = 82768.57000000001
$ws.Range("N138").Value = -93048.57000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 916.7
$ws.Range("I107").Value = 1013.5
$ws.Range("J107").Value = 529.5
$ws.Range("K107").Value = 3040.5
$ws.Range("L107").Value = 1588.5
$ws.Range("M107").Value = -1120.5
$ws.Range("N107").Value = -5428.5
$ws.Range("H122").Value = 17861010
$ws.Range("I122").Value = 3916.889
$ws.Range("J122").Value = 50003780
$ws.Range("K122").Value = 11750.667
$ws.Range("L122").Value = 150011340
$ws.Range("M122").Value = -9300.667000000001
$ws.Range("N122").Value = -150016240
$ws.Range("H132").Value = 10420031
$ws.Range("I132").Value = 12824142
$ws.Range("J132").Value = 2216.1667
$ws.Range("K132").Value = 38472426
$ws.Range("L132").Value = 6648.500100000001
$ws.Range("M132").Value = -38469896
$ws.Range("N132").Value = -11708.5001
$ws.Range("H136").Value = 21741354
$ws.Range("I136").Value = 27778824
$ws.Range("K136").Value = 83336472
$ws.Range("M136").Value = -83333922
